# Append: 2026-01-06 01:26 JST
# Update the "取得日時" (retrieved datetime) column for every data row
# on the "ランサーズ" sheet to the new timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2026-01-06 01:26:43"

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 1).Value = $newTimestamp
}
